$wb = $excel.ActiveWorkbook

# --- ALC row 8 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 93.40000000000001
$ws.Range("I8").Value = 93.40000000000001
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 280.2
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

# --- ALC row 17 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 529816.9
$ws.Range("J17").Value = 597978.4
$ws.Range("L17").Value = 1793935.2
$ws.Range("N17").Value = -1794271.2

# --- ALC row 33 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 6666815.5
$ws.Range("I33").Value = 7143005
$ws.Range("J33").Value = 159
$ws.Range("K33").Value = 7143005
$ws.Range("L33").Value = 159
$ws.Range("M33").Value = -7142776
$ws.Range("N33").Value = -617

# --- ALC row 39 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 344.33334
$ws.Range("I39").Value = 245.83333
$ws.Range("K39").Value = 737.49999
$ws.Range("M39").Value = -441.49999

# --- ALC row 40 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1106.862
$ws.Range("I40").Value = 1005
$ws.Range("K40").Value = 1005
$ws.Range("M40").Value = -830

# --- ALC row 42 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 100.42857
$ws.Range("I42").Value = 48.2
$ws.Range("K42").Value = 144.6
$ws.Range("M42").Value = 85.39999999999998

# --- ALC row 80 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 989.4524
$ws.Range("I80").Value = 1056.0555
$ws.Range("J80").Value = 939.5
$ws.Range("K80").Value = 3168.1665
$ws.Range("L80").Value = 2818.5
$ws.Range("M80").Value = -2170.1665
$ws.Range("N80").Value = -4814.5

# --- ALC row 83 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 989.4524
$ws.Range("I83").Value = 1056.0555
$ws.Range("J83").Value = 939.5
$ws.Range("K83").Value = 9504.4995
$ws.Range("L83").Value = 8455.5
$ws.Range("M83").Value = -4512.4995
$ws.Range("N83").Value = -18439.5

# --- ALC row 92 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 5239.476
$ws.Range("I92").Value = 5744.421
$ws.Range("K92").Value = 5744.421
$ws.Range("M92").Value = -4496.421

# --- ALC row 132 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 92648.17999999999
$ws.Range("I132").Value = 51813
$ws.Range("K132").Value = 155439
$ws.Range("M132").Value = -152909

# --- ALC row 137 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1700.8334
$ws.Range("I137").Value = 1675.2
$ws.Range("J137").Value = 1829
$ws.Range("K137").Value = 5025.6
$ws.Range("L137").Value = 5487
$ws.Range("M137").Value = -2475.6
$ws.Range("N137").Value = -10587

# --- ARM row 17 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 900
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 900
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -1246

# --- ARM row 32 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2316.3333
$ws.Range("I32").Value = 1914.2361
$ws.Range("K32").Value = 1914.2361
$ws.Range("M32").Value = -1627.2361

# --- ARM row 74 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5189.875
$ws.Range("I74").Value = 640.17645
$ws.Range("K74").Value = 640.17645
$ws.Range("M74").Value = 233.82355

# --- ARM row 77 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5189.875
$ws.Range("I77").Value = 640.17645
$ws.Range("K77").Value = 3200.88225
$ws.Range("M77").Value = 1167.11775

# --- ARM row 97 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2499.0833
$ws.Range("I97").Value = 2723.9
$ws.Range("J97").Value = 1375
$ws.Range("K97").Value = 2723.9
$ws.Range("L97").Value = 1375
$ws.Range("M97").Value = -2227.9
$ws.Range("N97").Value = -2367

# --- ARM row 102 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3083.818
$ws.Range("I102").Value = 3083.818
$ws.Range("K102").Value = 3083.818
$ws.Range("M102").Value = -1461.818

# --- ARM row 132 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2561.1904
$ws.Range("I132").Value = 2339.25
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 7017.75
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -4487.75
$ws.Range("N132").Value = -26060

# --- BSM row 86 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4396
$ws.Range("I86").Value = 2858.25
$ws.Range("J86").Value = 8496.666999999999
$ws.Range("K86").Value = 2858.25
$ws.Range("L86").Value = 8496.666999999999
$ws.Range("M86").Value = -1735.25
$ws.Range("N86").Value = -10742.667

# --- BSM row 89 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4396
$ws.Range("I89").Value = 2858.25
$ws.Range("J89").Value = 8496.666999999999
$ws.Range("K89").Value = 14291.25
$ws.Range("L89").Value = 42483.335
$ws.Range("M89").Value = -8675.25
$ws.Range("N89").Value = -53715.335

# --- BSM row 107 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 924.21875
$ws.Range("I107").Value = 919.2
$ws.Range("J107").Value = 999.5
$ws.Range("K107").Value = 919.2
$ws.Range("L107").Value = 999.5
$ws.Range("M107").Value = 1000.8
$ws.Range("N107").Value = -4839.5

# --- CRP row 31 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 50437.863
$ws.Range("I31").Value = 64827.062
$ws.Range("J31").Value = 12066.667
$ws.Range("K31").Value = 64827.062
$ws.Range("L31").Value = 12066.667
$ws.Range("M31").Value = -64532.062
$ws.Range("N31").Value = -12656.667

# --- CRP row 34 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 50437.863
$ws.Range("I34").Value = 64827.062
$ws.Range("J34").Value = 12066.667
$ws.Range("K34").Value = 64827.062
$ws.Range("L34").Value = 12066.667
$ws.Range("M34").Value = -64625.062
$ws.Range("N34").Value = -12470.667

# --- CRP row 42 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 23310.2
$ws.Range("I42").Value = 5278
$ws.Range("J42").Value = 35331.668
$ws.Range("K42").Value = 5278
$ws.Range("L42").Value = 35331.668
$ws.Range("M42").Value = -4685
$ws.Range("N42").Value = -36517.668

# --- CRP row 63 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 32166.666
$ws.Range("J63").Value = 46250
$ws.Range("L63").Value = 46250
$ws.Range("N63").Value = -47622

# --- CRP row 66 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H66").Value = 32166.666
$ws.Range("J66").Value = 46250
$ws.Range("L66").Value = 138750
$ws.Range("N66").Value = -145614

# --- CUL row 3 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6221.3335
$ws.Range("I3").Value = 1832
$ws.Range("K3").Value = 5496
$ws.Range("M3").Value = -5384

# --- CUL row 5 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6126.3335
$ws.Range("J5").Value = 17289.334
$ws.Range("L5").Value = 51868.00199999999
$ws.Range("N5").Value = -52092.00199999999

# --- CUL row 20 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

# --- CUL row 98 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 550.375
$ws.Range("J98").Value = 338.5
$ws.Range("L98").Value = 1015.5
$ws.Range("N98").Value = -4011.5

# --- CUL row 135 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 6126.3335
$ws.Range("J135").Value = 17289.334
$ws.Range("L135").Value = 155604.006
$ws.Range("N135").Value = -160674.006

# --- GSM row 80 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2381.875
$ws.Range("I80").Value = 1857
$ws.Range("J80").Value = 3256.6667
$ws.Range("K80").Value = 1857
$ws.Range("L80").Value = 3256.6667
$ws.Range("M80").Value = -859
$ws.Range("N80").Value = -5252.6667

# --- GSM row 83 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2381.875
$ws.Range("I83").Value = 1857
$ws.Range("J83").Value = 3256.6667
$ws.Range("K83").Value = 9285
$ws.Range("L83").Value = 16283.3335
$ws.Range("M83").Value = -4293
$ws.Range("N83").Value = -26267.3335

# --- LTW row 42 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 12259.643
$ws.Range("I42").Value = 8995
$ws.Range("J42").Value = 18136
$ws.Range("K42").Value = 8995
$ws.Range("L42").Value = 18136
$ws.Range("M42").Value = -8432
$ws.Range("N42").Value = -19262

# --- LTW row 49 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H49").Value = 12259.643
$ws.Range("I49").Value = 8995
$ws.Range("J49").Value = 18136
$ws.Range("K49").Value = 8995
$ws.Range("L49").Value = 18136
$ws.Range("M49").Value = -8848
$ws.Range("N49").Value = -18430

# --- LTW row 139 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 92000
$ws.Range("K139").Value = 92000
$ws.Range("N139").Value = -102280
